$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# ALC row 9
$ws_ALC.Range("H9").Value = 1569.6364
$ws_ALC.Range("I9").Value = 745
$ws_ALC.Range("K9").Value = 745
$ws_ALC.Range("M9").Value = -576

# ALC row 32
$ws_ALC.Range("H32").Value = 5427.857
$ws_ALC.Range("I32").Value = 998
$ws_ALC.Range("K32").Value = 998
$ws_ALC.Range("M32").Value = -672

# ALC row 57
$ws_ALC.Range("H57").Value = 44494.5
$ws_ALC.Range("I57").Value = 39999
$ws_ALC.Range("J57").Value = 48990
$ws_ALC.Range("K57").Value = 119997
$ws_ALC.Range("L57").Value = 146970
$ws_ALC.Range("M57").Value = -119498
$ws_ALC.Range("N57").Value = -147968

# ALC row 64
$ws_ALC.Range("H64").Value = 6917.778
$ws_ALC.Range("J64").Value = 7526.857
$ws_ALC.Range("L64").Value = 7526.857
$ws_ALC.Range("N64").Value = -8022.857

# ALC row 67
$ws_ALC.Range("H67").Value = 6917.778
$ws_ALC.Range("J67").Value = 7526.857
$ws_ALC.Range("L67").Value = 7526.857
$ws_ALC.Range("N67").Value = -9242.857

# ALC row 74
$ws_ALC.Range("H74").Value = 3628.1428
$ws_ALC.Range("I74").Value = 2375
$ws_ALC.Range("J74").Value = 5299
$ws_ALC.Range("K74").Value = 2375
$ws_ALC.Range("L74").Value = 5299
$ws_ALC.Range("M74").Value = -1439
$ws_ALC.Range("N74").Value = -7171

# ALC row 77
$ws_ALC.Range("H77").Value = 3628.1428
$ws_ALC.Range("I77").Value = 2375
$ws_ALC.Range("J77").Value = 5299
$ws_ALC.Range("K77").Value = 11875
$ws_ALC.Range("L77").Value = 26495
$ws_ALC.Range("M77").Value = -7195
$ws_ALC.Range("N77").Value = -35855

# ALC row 86
$ws_ALC.Range("H86").Value = 8945.2
$ws_ALC.Range("I86").Value = 13752.286
$ws_ALC.Range("J86").Value = 4739
$ws_ALC.Range("K86").Value = 13752.286
$ws_ALC.Range("L86").Value = 4739
$ws_ALC.Range("M86").Value = -12629.286
$ws_ALC.Range("N86").Value = -6985

# ALC row 89
$ws_ALC.Range("H89").Value = 8945.2
$ws_ALC.Range("I89").Value = 13752.286
$ws_ALC.Range("J89").Value = 4739
$ws_ALC.Range("K89").Value = 68761.43
$ws_ALC.Range("L89").Value = 23695
$ws_ALC.Range("M89").Value = -63145.42999999999
$ws_ALC.Range("N89").Value = -34927

$ws_ARM = $wb.Worksheets.Item("ARM")
# ARM row 63
$ws_ARM.Range("H63").Value = 1640.8334
$ws_ARM.Range("I63").Value = 1587.25
$ws_ARM.Range("K63").Value = 1587.25
$ws_ARM.Range("M63").Value = -901.25

# ARM row 66
$ws_ARM.Range("H66").Value = 1640.8334
$ws_ARM.Range("I66").Value = 1587.25
$ws_ARM.Range("K66").Value = 7936.25
$ws_ARM.Range("M66").Value = -4504.25

$ws_BSM = $wb.Worksheets.Item("BSM")
# BSM row 22
$ws_BSM.Range("H22").Value = 15885401
$ws_BSM.Range("I22").Value = 17871050
$ws_BSM.Range("K22").Value = 17871050
$ws_BSM.Range("M22").Value = -17870877

# BSM row 86
$ws_BSM.Range("H86").Value = 50002540
$ws_BSM.Range("I86").Value = 2769
$ws_BSM.Range("J86").Value = 333334600
$ws_BSM.Range("K86").Value = 2769
$ws_BSM.Range("L86").Value = 333334600
$ws_BSM.Range("M86").Value = -1646
$ws_BSM.Range("N86").Value = -333336846

# BSM row 89
$ws_BSM.Range("H89").Value = 50002540
$ws_BSM.Range("I89").Value = 2769
$ws_BSM.Range("J89").Value = 333334600
$ws_BSM.Range("K89").Value = 13845
$ws_BSM.Range("L89").Value = 1666673000
$ws_BSM.Range("M89").Value = -8229
$ws_BSM.Range("N89").Value = -1666684232

$ws_CRP = $wb.Worksheets.Item("CRP")
# CRP row 19
$ws_CRP.Range("H19").Value = 220.25
$ws_CRP.Range("I19").Value = 220.25
$ws_CRP.Range("K19").Value = 220.25
$ws_CRP.Range("M19").Value = -50.25

# CRP row 24
$ws_CRP.Range("H24").Value = 220.25
$ws_CRP.Range("I24").Value = 220.25
$ws_CRP.Range("K24").Value = 220.25
$ws_CRP.Range("M24").Value = -50.25

# CRP row 31
$ws_CRP.Range("H31").Value = 4815.1665
$ws_CRP.Range("I31").Value = 1021.76086
$ws_CRP.Range("J31").Value = 26627.25
$ws_CRP.Range("K31").Value = 1021.76086
$ws_CRP.Range("L31").Value = 26627.25
$ws_CRP.Range("M31").Value = -726.76086
$ws_CRP.Range("N31").Value = -27217.25

# CRP row 34
$ws_CRP.Range("H34").Value = 4815.1665
$ws_CRP.Range("I34").Value = 1021.76086
$ws_CRP.Range("J34").Value = 26627.25
$ws_CRP.Range("K34").Value = 1021.76086
$ws_CRP.Range("L34").Value = 26627.25
$ws_CRP.Range("M34").Value = -819.76086
$ws_CRP.Range("N34").Value = -27031.25

# CRP row 62
$ws_CRP.Range("H62").Value = 7000
$ws_CRP.Range("I62").Value = 6000
$ws_CRP.Range("K62").Value = 6000
$ws_CRP.Range("M62").Value = -5376

# CRP row 65
$ws_CRP.Range("H65").Value = 7000
$ws_CRP.Range("I65").Value = 6000
$ws_CRP.Range("K65").Value = 30000
$ws_CRP.Range("M65").Value = -26880

# CRP row 132
$ws_CRP.Range("H132").Value = 34484884
$ws_CRP.Range("I132").Value = 2208.3704
$ws_CRP.Range("K132").Value = 6625.111199999999
$ws_CRP.Range("M132").Value = -4095.111199999999

$ws_CUL = $wb.Worksheets.Item("CUL")
# CUL row 14
$ws_CUL.Range("H14").Value = 67.5
$ws_CUL.Range("I14").Value = 67.5
$ws_CUL.Range("K14").Value = 202.5
$ws_CUL.Range("M14").Value = -29.5

# CUL row 103
$ws_CUL.Range("H103").Value = 843.25
$ws_CUL.Range("J103").Value = 457.66666
$ws_CUL.Range("L103").Value = 1372.99998
$ws_CUL.Range("N103").Value = -3130.99998

# CUL row 124
$ws_CUL.Range("H124").Value = 5680.8
$ws_CUL.Range("I124").Value = 5680.8
$ws_CUL.Range("J124").Value = 0
$ws_CUL.Range("K124").Value = 17042.4
$ws_CUL.Range("L124").Value = 0
$ws_CUL.Range("M124").Value = -12132.4
$ws_CUL.Range("N124").ClearContents()

# CUL row 129
$ws_CUL.Range("H129").Value = 6062429.5
$ws_CUL.Range("J129").Value = 18183902
$ws_CUL.Range("L129").Value = 54551706
$ws_CUL.Range("N129").Value = -54561706

$ws_GSM = $wb.Worksheets.Item("GSM")
# GSM row 126
$ws_GSM.Range("H126").Value = 5371112.5
$ws_GSM.Range("I126").Value = 2569209.8
$ws_GSM.Range("K126").Value = 7707629.399999999
$ws_GSM.Range("M126").Value = -7705159.399999999

# GSM row 140
$ws_GSM.Range("H140").Value = 0
$ws_GSM.Range("J140").Value = 0
$ws_GSM.Range("L140").Value = 0
$ws_GSM.Range("N140").ClearContents()

$ws_LTW = $wb.Worksheets.Item("LTW")
# LTW row 16
$ws_LTW.Range("H16").Value = 100001800
$ws_LTW.Range("I16").Value = 111112940
$ws_LTW.Range("J16").Value = 1500
$ws_LTW.Range("K16").Value = 111112940
$ws_LTW.Range("L16").Value = 1500
$ws_LTW.Range("M16").Value = -111112770
$ws_LTW.Range("N16").Value = -1840

# LTW row 40
$ws_LTW.Range("H40").Value = 5352093.5
$ws_LTW.Range("I40").Value = 1000
$ws_LTW.Range("J40").Value = 7358753.5
$ws_LTW.Range("K40").Value = 1000
$ws_LTW.Range("L40").Value = 7358753.5
$ws_LTW.Range("M40").Value = -864
$ws_LTW.Range("N40").Value = -7359025.5

# LTW row 99
$ws_LTW.Range("H99").Value = 15259
$ws_LTW.Range("I99").Value = 15259
$ws_LTW.Range("J99").Value = 0
$ws_LTW.Range("K99").Value = 15259
$ws_LTW.Range("L99").Value = 0
$ws_LTW.Range("N99").ClearContents()
$ws_LTW.Range("M99").Value = -12264

# LTW row 132
$ws_LTW.Range("H132").Value = 1546296.2
$ws_LTW.Range("I132").Value = 3737.2
$ws_LTW.Range("J132").Value = 6688160
$ws_LTW.Range("K132").Value = 11211.6
$ws_LTW.Range("L132").Value = 20064480
$ws_LTW.Range("M132").Value = -8681.599999999999
$ws_LTW.Range("N132").Value = -20069540

$ws_WVR = $wb.Worksheets.Item("WVR")
# WVR row 81
$ws_WVR.Range("H81").Value = 5000
$ws_WVR.Range("I81").Value = 0
$ws_WVR.Range("J81").Value = 5000
$ws_WVR.Range("K81").Value = 0
$ws_WVR.Range("L81").Value = 10000
$ws_WVR.Range("M81").ClearContents()
$ws_WVR.Range("N81").Value = -12122

# WVR row 84
$ws_WVR.Range("H84").Value = 5000
$ws_WVR.Range("I84").Value = 0
$ws_WVR.Range("J84").Value = 5000
$ws_WVR.Range("K84").Value = 0
$ws_WVR.Range("L84").Value = 50000
$ws_WVR.Range("M84").ClearContents()
$ws_WVR.Range("N84").Value = -60608
